$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per diff
$ws.Range("D2").Value = 0.002571262486839021
$ws.Range("E2").Value = 0.09573784751843535

$ws.Range("D3").Value = 0.00195636559570375

$ws.Range("D4").Value = 0.005938513049319022

$ws.Range("D5").Value = 0.007707478666925026

# Add new row 6 (DWA)
$ws.Range("A6").Value = "DWA"
$ws.Range("B6").Value = 45
$ws.Range("C6").Value = 0.45
$ws.Range("D6").Value = 0.008751289595436359
$ws.Range("E6").Value = 0.4078261088489451

# Apply same style as the other rows' A column cells (A2:A5) to A6
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
